$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("L2").Value = 7.5
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.67
$ws.Range("AO2").Value = 7.5
$ws.Range("BA2").Value = 251

# Row 3 updates
$ws.Range("G3").Value = 1.5
$ws.Range("I3").Value = 6.5
$ws.Range("L3").Value = 7
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 2.63
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
$ws.Range("W3").Value = 5.5
$ws.Range("Z3").Value = 10
$ws.Range("AC3").Value = 8
$ws.Range("AE3").Value = 23
$ws.Range("AF3").Value = 81
$ws.Range("AJ3").Value = 21
$ws.Range("AL3").Value = 51
$ws.Range("AN3").Value = 3.25
$ws.Range("AQ3").Value = 26
$ws.Range("AS3").Value = 201
$ws.Range("AT3").Value = 2.63
$ws.Range("AY3").Value = 41
$ws.Range("AZ3").Value = 151
$ws.Range("BA3").Value = 201
